# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 198
    "F3"  = 434
    "F4"  = 12738
    "F5"  = 1330
    "F6"  = 168
    "F9"  = 164
    "F10" = 213
    "F11" = 460
    "F12" = 62
    "F16" = 388
    "F17" = 5456
    "F18" = 101
    "F19" = 30
    "F20" = 950
    "F21" = 26
    "F23" = 96
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
